$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing 200m - M results
$ws.Range("B4").Value = "Letsile Tebogo"
$ws.Range("C4").Value = "Kenneth Bednarek"
$ws.Range("D4").Value = "Noah Lyles"

# Move the active selection to D5 (matches the saved view state)
$ws.Range("D5").Select()
